$wb = $excel.ActiveWorkbook

# --- Update raw input data (columns D:I, rows 2-7) on the "EtOAc Ethanol Water" sheet ---
# Dependent formula columns (J:U) recalc automatically.
$wsEtOAc = $wb.Worksheets.Item("EtOAc Ethanol Water")
$wsEtOAc.Activate()

$data = New-Object 'object[,]' 6,6
$data[0,0] = 0.78719393
$data[0,1] = 0
$data[0,2] = 0.21280607
$data[0,3] = 0.01890998
$data[0,4] = 0
$data[0,5] = 0.98109002
$data[1,0] = 0.70216729
$data[1,1] = 0.0524376
$data[1,2] = 0.2453951
$data[1,3] = 0.02341177
$data[1,4] = 0.02032564
$data[1,5] = 0.95626259
$data[2,0] = 0.62301631
$data[2,1] = 0.09934726
$data[2,2] = 0.27763643
$data[2,3] = 0.02888137
$data[2,4] = 0.04204784
$data[2,5] = 0.92907078
$data[3,0] = 0.55059846
$data[3,1] = 0.14064636
$data[3,2] = 0.30875518
$data[3,3] = 0.03557354
$data[3,4] = 0.06554427
$data[3,5] = 0.89888219
$data[4,0] = 0.47564347
$data[4,1] = 0.17805133
$data[4,2] = 0.3463052
$data[4,3] = 0.04280507
$data[4,4] = 0.09053062
$data[4,5] = 0.86666431
$data[5,0] = 0.36121251
$data[5,1] = 0.21489175
$data[5,2] = 0.42389574
$data[5,3] = 0.04862695
$data[5,4] = 0.11571994
$data[5,5] = 0.83565311

$wsEtOAc.Range("D2:I7").Value = $data

# --- Update sheet selection/activation state ---
# Before: "D-Limonene Ethanol Water" tab was active/selected.
# After: "EtOAc Ethanol Water" tab is active/selected, with a new cell selection.
$wsEtOAc.Range("I11").Select()

Write-Host "Edit applied"
